# Add data for 2021-12-21 (carjacking-by-neighborhood-by-month.xlsx)
# - rename sheet / update "through" header to reflect the new as-of date
# - add a few newly-populated cells
# - bump several existing December year-over-year counters by 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet title + column header text: "through December 20" -> "through December 21"
$ws.Name = "Through 2021-12-21"
$ws.Range("B1").Value = "December 2021 (through December 21)"

# Newly populated cells (previously blank)
$ws.Range("B2").Value = 1
$ws.Range("Z13").Value = 1
$ws.Range("BJ18").Value = 1
$ws.Range("Z24").Value = 2
$ws.Range("N46").Value = 1
$ws.Range("AL95").Value = 1

# Updated counts on existing cells
$ws.Range("N4").Value = 10
$ws.Range("Z4").Value = 5
$ws.Range("Z6").Value = 3
$ws.Range("BV6").Value = 2
$ws.Range("AL7").Value = 10
$ws.Range("AX7").Value = 9
$ws.Range("Z11").Value = 2
$ws.Range("AL14").Value = 2
$ws.Range("N18").Value = 4
$ws.Range("BJ30").Value = 3
$ws.Range("N38").Value = 2
$ws.Range("AX41").Value = 2
$ws.Range("N61").Value = 2
